$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 271.25
$ws.Range("I6").Value = 271.25
$ws.Range("K6").Value = 813.75
$ws.Range("M6").Value = -701.75

$ws.Range("H111").Value = 3897.375
$ws.Range("J111").Value = 4450
$ws.Range("L111").Value = 13350
$ws.Range("N111").Value = -19484

$ws.Range("H132").Value = 419725.56
$ws.Range("I132").Value = 450686.3
$ws.Range("K132").Value = 1352058.9
$ws.Range("M132").Value = -1349528.9

$ws.Range("H137").Value = 3303.1667
$ws.Range("I137").Value = 2220.625
$ws.Range("J137").Value = 3844.4375
$ws.Range("K137").Value = 6661.875
$ws.Range("L137").Value = 11533.3125
$ws.Range("M137").Value = -4111.875
$ws.Range("N137").Value = -16633.3125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6179038.5
$ws.Range("I32").Value = 7147019
$ws.Range("K32").Value = 7147019
$ws.Range("M32").Value = -7146732

$ws.Range("H45").Value = 3126.2
$ws.Range("I45").Value = 2916.4546
$ws.Range("K45").Value = 2916.4546
$ws.Range("M45").Value = -2539.4546

$ws.Range("H125").Value = 99500
$ws.Range("J125").Value = 99500
$ws.Range("L125").Value = 99500
$ws.Range("N125").Value = -109340

$ws.Range("H140").Value = 89000
$ws.Range("J140").Value = 89000
$ws.Range("L140").Value = 89000
$ws.Range("N140").Value = -99360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").ClearContents()
$ws.Range("N9").Value = 0

$ws.Range("H47").Value = 755000
$ws.Range("J47").Value = 755000
$ws.Range("L47").Value = 755000
$ws.Range("N47").Value = -756040

$ws.Range("H140").Value = 105666.11
$ws.Range("I140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("M140").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 16029.059
$ws.Range("I6").Value = 10549.4
$ws.Range("K6").Value = 10549.4
$ws.Range("M6").Value = -10436.4

$ws.Range("H31").Value = 7239.6206
$ws.Range("I31").Value = 1442.5454
$ws.Range("J31").Value = 10782.277
$ws.Range("K31").Value = 1442.5454
$ws.Range("L31").Value = 10782.277
$ws.Range("M31").Value = -1147.5454
$ws.Range("N31").Value = -11372.277

$ws.Range("H32").Value = 5000
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()

$ws.Range("H34").Value = 7239.6206
$ws.Range("I34").Value = 1442.5454
$ws.Range("J34").Value = 10782.277
$ws.Range("K34").Value = 1442.5454
$ws.Range("L34").Value = 10782.277
$ws.Range("M34").Value = -1240.5454
$ws.Range("N34").Value = -11186.277

$ws.Range("H134").Value = 38649.77
$ws.Range("I134").Value = 31870.584
$ws.Range("K134").Value = 95611.75199999999
$ws.Range("M134").Value = -93076.75199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 5857907
$ws.Range("I7").Value = 6666715.5
$ws.Range("K7").Value = 20000146.5
$ws.Range("M7").Value = -20000034.5

$ws.Range("H80").Value = 5495
$ws.Range("I80").Value = 4995
$ws.Range("J80").Value = 5995
$ws.Range("K80").Value = 14985
$ws.Range("L80").Value = 17985
$ws.Range("M80").Value = -14049
$ws.Range("N80").Value = -19857

$ws.Range("H81").Value = 6212.1665
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 6212.1665
$ws.Range("K81").Value = 0
$ws.Range("L81").ClearContents()
$ws.Range("M81").Value = 18636.4995
$ws.Range("N81").Value = -20882.4995

$ws.Range("H83").Value = 5495
$ws.Range("I83").Value = 4995
$ws.Range("J83").Value = 5995
$ws.Range("K83").Value = 44955
$ws.Range("L83").Value = 53955
$ws.Range("M83").Value = -40275
$ws.Range("N83").Value = -63315

$ws.Range("H84").Value = 6212.1665
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 6212.1665
$ws.Range("K84").Value = 0
$ws.Range("L84").ClearContents()
$ws.Range("M84").Value = 55909.4985
$ws.Range("N84").Value = -67141.4985

$ws.Range("H86").Value = 859.1818
$ws.Range("J86").Value = 750.1429000000001
$ws.Range("L86").Value = 2250.4287
$ws.Range("N86").Value = -4622.4287

$ws.Range("H89").Value = 859.1818
$ws.Range("J89").Value = 750.1429000000001
$ws.Range("L89").Value = 6751.2861
$ws.Range("N89").Value = -18607.2861

$ws.Range("H92").Value = 805.0769
$ws.Range("I92").Value = 830.25
$ws.Range("J92").Value = 503
$ws.Range("K92").Value = 2490.75
$ws.Range("L92").Value = 1509
$ws.Range("M92").Value = -1242.75
$ws.Range("N92").Value = -4005

$ws.Range("H107").Value = 995.2222
$ws.Range("J107").Value = 1199.8
$ws.Range("L107").Value = 3599.4
$ws.Range("N107").Value = -7439.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").ClearContents()
$ws.Range("N75").Value = 0

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").ClearContents()
$ws.Range("N78").Value = 0

$ws.Range("H80").Value = 166422.64
$ws.Range("I80").Value = 266921.78
$ws.Range("J80").Value = 7299
$ws.Range("K80").Value = 266921.78
$ws.Range("L80").Value = 7299
$ws.Range("M80").Value = -265923.78
$ws.Range("N80").Value = -9295

$ws.Range("H83").Value = 166422.64
$ws.Range("I83").Value = 266921.78
$ws.Range("J83").Value = 7299
$ws.Range("K83").Value = 1334608.9
$ws.Range("L83").Value = 36495
$ws.Range("M83").Value = -1329616.9
$ws.Range("N83").Value = -46479

$ws.Range("H102").Value = 2384.5557
$ws.Range("I102").Value = 1933.0834
$ws.Range("K102").Value = 1933.0834
$ws.Range("M102").Value = -311.0834

$ws.Range("H123").Value = 80000
$ws.Range("J123").Value = 80000
$ws.Range("L123").Value = 80000
$ws.Range("N123").Value = -84900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H138").Value = 120000
$ws.Range("J138").Value = 120000
$ws.Range("L138").Value = 120000
$ws.Range("N138").Value = -130280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 11904.75
$ws.Range("J45").Value = 11904.75
$ws.Range("L45").Value = 11904.75
$ws.Range("N45").Value = -12886.75

$ws.Range("H86").Value = 78999.75
$ws.Range("J86").Value = 78999.75
$ws.Range("L86").Value = 78999.75
$ws.Range("N86").Value = -81245.75

$ws.Range("H89").Value = 78999.75
$ws.Range("J89").Value = 78999.75
$ws.Range("L89").Value = 394998.75
$ws.Range("N89").Value = -406230.75
